# Reorder/refresh the resource summary rows (rows 3-14) on the "Resources"
# sheet to reflect the latest values read from the tracking periods parse.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resources")

# Target state for columns A:H, rows 3 through 14 (ID, Name, Type,
# Availability, Cost/Use, Cost/Unit, Assigned To, Total Cost)
$rows = @(
    @(12, "joiner",              "RENEWABLE", "5 #5", 0, 42,    "40[5 #5];", 8400),
    @(6,  "roofer",               "RENEWABLE", "4 #4", 0, 36,    "49[4 #4];", 11520),
    @(1,  "plumber",              "RENEWABLE", "4 #4", 0, 36,    "30;52;53;", 8640),
    @(3,  "team subcontractor",   "RENEWABLE", "8 #8", 0, 38.56, "6[4 #8];17;18[3 #8];19[4 #8];20[4 #8];21[3 #8];22[3 #8];23[5 #8];24[4 #8];25[2 #8];28[5 #8];69;", 173057.28),
    @(4,  "mason",                "RENEWABLE", "4 #4", 0, 40.65, "27;45;55;", 9756),
    @(2,  "electrician",          "RENEWABLE", "4 #4", 0, 32,    "29;51;54;61;", 24320),
    @(5,  "heating engineer",     "RENEWABLE", "2 #2", 0, 43.2,  "59;60;", 13824),
    @(9,  "screed layer",         "RENEWABLE", "3 #3", 0, 0,     "34[2 #3];56;", 0),
    @(7,  "carpenter",            "RENEWABLE", "3 #3", 0, 38,    "26[3 #3];", 9120),
    @(11, "painter",              "RENEWABLE", "6 #6", 0, 35,    "39[4 #6];", 16800),
    @(8,  "floor layer",          "RENEWABLE", "4 #4", 0, 41.56, "35;57;58;", 19948.8),
    @(10, "plasterer",            "RENEWABLE", "4 #4", 0, 41.38, "32[3 #4];", 14896.8)
)

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
